# Apply the "Add data for 2022-08-03" update to the carjacking-by-neighborhood
# workbook: rename the sheet / header label from "July 25" to "July 26" and
# update the July-2022-to-date counts for several neighborhoods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet tab and the matching column header label ---
$ws.Name = "Through 2022-07-26"
$ws.Range("B1").Value = "July 2022 (through July 26)"

# --- Update / add neighborhood counts for the "July 2022 (through ..)" column (B) ---
$ws.Range("B2").Value = 13
$ws.Range("B7").Value = 4
$ws.Range("B33").Value = 2
$ws.Range("B38").Value = 4
$ws.Range("B52").Value = 5
$ws.Range("B57").Value = 3

# --- Other monthly cell updates across the sheet ---
$ws.Range("AY3").Value = 3
$ws.Range("W5").Value = 3
$ws.Range("AK6").Value = 1
$ws.Range("P7").Value = 1
$ws.Range("AY7").Value = 2
$ws.Range("I14").Value = 3
$ws.Range("P18").Value = 1
$ws.Range("AK19").Value = 3
$ws.Range("AK29").Value = 1
$ws.Range("P32").Value = 1
$ws.Range("AK41").Value = 3
$ws.Range("P45").Value = 3
$ws.Range("AY57").Value = 1
$ws.Range("I65").Value = 3
